# The "P53" sheet had a cell style (numFmtId 0 / applyNumberFormat) applied
# across most of its data range; it was never actually needed, so clear the
# formatting back to the default style -- but do it range-by-range (rather
# than as one C2:L8 rectangle) so we don't touch the always-empty columns
# G:H / rows 5:6 that sit inside that bounding box but were never part of
# the sheet's real data.
$wb = $excel.ActiveWorkbook
$sheet = $wb.Worksheets.Item("P53")

$formattedRanges = @(
    "C2:F2", "I2:L2",
    "C3:D3", "I3:L3",
    "C7:F7", "I7:L7",
    "C8:F8", "I8:L8"
)
foreach ($addr in $formattedRanges) {
    $sheet.Range($addr).ClearFormats()
}

# E3 and F3 only ever held an applied style with no value -- drop them
# entirely rather than leaving blank, still-styled cells behind.
$sheet.Range("E3:F3").Clear()

# The sheet was mislabeled "P53"; correct the casing to match the other
# p53-related assets.
$sheet.Name = "p53"

# Leave the cursor where the author last left it when they saved.
$sheet.Range("F14").Select()
